$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AE5").Value = 0.4
$ws.Range("AF5").Value = 0.3
$ws.Range("AG5").Value = 0.2
$ws.Range("AH5").Value = 0.8

$ws.Range("AH8").Select()
